# "add current stock price to company data"
#
# The "companies" table (Sheet1, starting around row 59) gets two new
# COLUMN / DATA TYPE rows inserted right before the existing
# "3_month_trading_volume" row:
#   current_price         | string
#   previous_close_price  | string
#
# Inserting the two rows pushes every row below them (everything from the
# old row 71 onward) down by two, which matches the rest of the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the "3_month_trading_volume" row (row 71).
$ws.Rows.Item(71).Resize(2).Insert()

# The row that used to be 71 ("3_month_trading_volume") is now row 73 - copy
# its B:D formatting (borders/styles) into the two freshly-inserted rows so
# they look like the rest of the table.
$ws.Range("B73:D73").Copy()
$ws.Range("B71:D72").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new column / data-type cells.
$ws.Range("B71").Value = "current_price"
$ws.Range("C71").Value = "string"

$ws.Range("B72").Value = "previous_close_price"
$ws.Range("C72").Value = "string"

# Reflect the author's saved selection/viewport.
$ws.Range("C72").Select()
